$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 150.0354306666667
$ws.Range("H2").Value = 450.106292
$ws.Range("I2").Value = 0.4152507364956075
$ws.Range("J2").Value = 0.4152507364956075
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 833.4308676666666
$ws.Range("N2").Value = 2500.292603
$ws.Range("O2").Value = 0.8518935545813505
$ws.Range("P2").Value = 0.8518935545813505
$ws.Range("Q2").Value = 125044.159161262
$ws.Range("R2").Value = 1125397.432451358
$ws.Range("S2").Value = 0.3537494259557668
$ws.Range("T2").Value = 0.3537494259557667

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 150.0354306666667
$ws.Range("H3").Value = 450.106292
$ws.Range("I3").Value = 0.4152507364956075
$ws.Range("J3").Value = 0.4152507364956075
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 17.73945766666667
$ws.Range("N3").Value = 53.218373
$ws.Range("O3").Value = 0.01813243333584592
$ws.Range("P3").Value = 0.01813243333584592
$ws.Range("Q3").Value = 2661.547170811435
$ws.Range("R3").Value = 23953.92453730292
$ws.Range("S3").Value = 0.007529506297167524
$ws.Range("T3").Value = 0.007529506297167522

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 150.0354306666667
$ws.Range("H4").Value = 450.106292
$ws.Range("I4").Value = 0.4152507364956075
$ws.Range("J4").Value = 0.4152507364956075
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.344749666666667
$ws.Range("N4").Value = 4.034249
$ws.Range("O4").Value = 0.001374539410528448
$ws.Range("P4").Value = 0.001374539410528448
$ws.Range("Q4").Value = 201.7600953771898
$ws.Range("R4").Value = 1815.840858394708
$ws.Range("S4").Value = 0.0005707785025641763
$ws.Range("T4").Value = 0.0005707785025641761

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 150.0354306666667
$ws.Range("H5").Value = 450.106292
$ws.Range("I5").Value = 0.4152507364956075
$ws.Range("J5").Value = 0.4152507364956075
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 125.812397
$ws.Range("N5").Value = 377.437191
$ws.Range("O5").Value = 0.1285994726722751
$ws.Range("P5").Value = 0.1285994726722751
$ws.Range("Q5").Value = 18876.31716710064
$ws.Range("R5").Value = 169886.8545039058
$ws.Range("S5").Value = 0.05340102574010901
$ws.Range("T5").Value = 0.053401025740109

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 68.382243
$ws.Range("H6").Value = 205.146729
$ws.Range("I6").Value = 0.1892604742946246
$ws.Range("J6").Value = 0.1892604742946246
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 833.4308676666666
$ws.Range("N6").Value = 2500.292603
$ws.Range("O6").Value = 0.8518935545813505
$ws.Range("P6").Value = 0.8518935545813505
$ws.Range("Q6").Value = 56991.87211648284
$ws.Range("R6").Value = 512926.8490483455
$ws.Range("S6").Value = 0.1612297781886001
$ws.Range("T6").Value = 0.1612297781886001

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 68.382243
$ws.Range("H7").Value = 205.146729
$ws.Range("I7").Value = 0.1892604742946246
$ws.Range("J7").Value = 0.1892604742946246
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 17.73945766666667
$ws.Range("N7").Value = 53.218373
$ws.Range("O7").Value = 0.01813243333584592
$ws.Range("P7").Value = 0.01813243333584592
$ws.Range("Q7").Value = 1213.063904850213
$ws.Range("R7").Value = 10917.57514365192
$ws.Range("S7").Value = 0.003431752933257862
$ws.Range("T7").Value = 0.003431752933257861

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 68.382243
$ws.Range("H8").Value = 205.146729
$ws.Range("I8").Value = 0.1892604742946246
$ws.Range("J8").Value = 0.1892604742946246
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.344749666666667
$ws.Range("N8").Value = 4.034249
$ws.Range("O8").Value = 0.001374539410528448
$ws.Range("P8").Value = 0.001374539410528448
$ws.Range("Q8").Value = 91.95699848016901
$ws.Range("R8").Value = 827.612986321521
$ws.Range("S8").Value = 0.0002601459807732679
$ws.Range("T8").Value = 0.0002601459807732678

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 68.382243
$ws.Range("H9").Value = 205.146729
$ws.Range("I9").Value = 0.1892604742946246
$ws.Range("J9").Value = 0.1892604742946246
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 125.812397
$ws.Range("N9").Value = 377.437191
$ws.Range("O9").Value = 0.1285994726722751
$ws.Range("P9").Value = 0.1285994726722751
$ws.Range("Q9").Value = 8603.333904066472
$ws.Range("R9").Value = 77430.00513659825
$ws.Range("S9").Value = 0.02433879719199341
$ws.Range("T9").Value = 0.02433879719199341

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 104.737245
$ws.Range("H10").Value = 314.211735
$ws.Range("I10").Value = 0.2898796499701289
$ws.Range("J10").Value = 0.2898796499701289
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 833.4308676666666
$ws.Range("N10").Value = 2500.292603
$ws.Range("O10").Value = 0.8518935545813505
$ws.Range("P10").Value = 0.8518935545813505
$ws.Range("Q10").Value = 87291.25297736622
$ws.Range("R10").Value = 785621.2767962961
$ws.Range("S10").Value = 0.2469466054138508
$ws.Range("T10").Value = 0.2469466054138508

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 104.737245
$ws.Range("H11").Value = 314.211735
$ws.Range("I11").Value = 0.2898796499701289
$ws.Range("J11").Value = 0.2898796499701289
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 17.73945766666667
$ws.Range("N11").Value = 53.218373
$ws.Range("O11").Value = 0.01813243333584592
$ws.Range("P11").Value = 0.01813243333584592
$ws.Range("Q11").Value = 1857.981923800795
$ws.Range("R11").Value = 16721.83731420715
$ws.Range("S11").Value = 0.005256223428501713
$ws.Range("T11").Value = 0.005256223428501712

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 104.737245
$ws.Range("H12").Value = 314.211735
$ws.Range("I12").Value = 0.2898796499701289
$ws.Range("J12").Value = 0.2898796499701289
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.344749666666667
$ws.Range("N12").Value = 4.034249
$ws.Range("O12").Value = 0.001374539410528448
$ws.Range("P12").Value = 0.001374539410528448
$ws.Range("Q12").Value = 140.845375301335
$ws.Range("R12").Value = 1267.608377712015
$ws.Range("S12").Value = 0.0003984510031941339
$ws.Range("T12").Value = 0.0003984510031941338

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 104.737245
$ws.Range("H13").Value = 314.211735
$ws.Range("I13").Value = 0.2898796499701289
$ws.Range("J13").Value = 0.2898796499701289
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 125.812397
$ws.Range("N13").Value = 377.437191
$ws.Range("O13").Value = 0.1285994726722751
$ws.Range("P13").Value = 0.1285994726722751
$ws.Range("Q13").Value = 13177.24384862627
$ws.Range("R13").Value = 118595.1946376364
$ws.Range("S13").Value = 0.03727837012458228
$ws.Range("T13").Value = 0.03727837012458227

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 38.15794
$ws.Range("H14").Value = 114.47382
$ws.Range("I14").Value = 0.105609139239639
$ws.Range("J14").Value = 0.105609139239639
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 833.4308676666666
$ws.Range("N14").Value = 2500.292603
$ws.Range("O14").Value = 0.8518935545813505
$ws.Range("P14").Value = 0.8518935545813505
$ws.Range("Q14").Value = 31802.00504257261
$ws.Range("R14").Value = 286218.0453831534
$ws.Range("S14").Value = 0.08996774502313286
$ws.Range("T14").Value = 0.08996774502313283

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 38.15794
$ws.Range("H15").Value = 114.47382
$ws.Range("I15").Value = 0.105609139239639
$ws.Range("J15").Value = 0.105609139239639
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 17.73945766666667
$ws.Range("N15").Value = 53.218373
$ws.Range("O15").Value = 0.01813243333584592
$ws.Range("P15").Value = 0.01813243333584592
$ws.Range("Q15").Value = 676.9011612772067
$ws.Range("R15").Value = 6092.11045149486
$ws.Range("S15").Value = 0.001914950676918824
$ws.Range("T15").Value = 0.001914950676918823

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 38.15794
$ws.Range("H16").Value = 114.47382
$ws.Range("I16").Value = 0.105609139239639
$ws.Range("J16").Value = 0.105609139239639
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.344749666666667
$ws.Range("N16").Value = 4.034249
$ws.Range("O16").Value = 0.001374539410528448
$ws.Range("P16").Value = 0.001374539410528448
$ws.Range("Q16").Value = 51.31287709568667
$ws.Range("R16").Value = 461.81589386118
$ws.Range("S16").Value = 0.0001451639239968702
$ws.Range("T16").Value = 0.0001451639239968701

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 38.15794
$ws.Range("H17").Value = 114.47382
$ws.Range("I17").Value = 0.105609139239639
$ws.Range("J17").Value = 0.105609139239639
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 125.812397
$ws.Range("N17").Value = 377.437191
$ws.Range("O17").Value = 0.1285994726722751
$ws.Range("P17").Value = 0.1285994726722751
$ws.Range("Q17").Value = 4800.741895982181
$ws.Range("R17").Value = 43206.67706383963
$ws.Range("S17").Value = 0.01358127961559046
$ws.Range("T17").Value = 0.01358127961559045

